# Rename the "_old"/"_new" suffixed column headers to the format-version
# specific suffixes "_FV2304" / "_FV2310", freeze the header row, and wrap
# the used range in an Excel Table ("Table1") with an AutoFilter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A-J: "<Name>_old" -> "<Name>_FV2304"
for ($i = 0; $i -lt $oldNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $oldNames[$i] + "_FV2304"
}

# Column K ("diff") is unchanged.

# Columns L-U: "<Name>_new" -> "<Name>_FV2310"
for ($i = 0; $i -lt $oldNames.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $oldNames[$i] + "_FV2310"
}

# Freeze the header row (row 1).
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)

# Turn the used range into a proper Excel Table with an AutoFilter.
$dataRange = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
[void]($tbl.TableStyle = "")
